# Updates the cryptos sheet: refreshed prices/volumes, the "Hora" (hour)
# column bumped from 17 to 18, and several coin rows re-synced to reflect
# the latest coinranking.com ordering (rows 11-20 and 41-43 shifted).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.01"
$ws.Range("G2").Value = "'18"
$ws.Range("D3").Value = "'24.97"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'5.050"
$ws.Range("G4").Value = "'18"
$ws.Range("D5").Value = "'0.05613"
$ws.Range("G5").Value = "'18"
$ws.Range("D6").Value = "'6.489"
$ws.Range("G6").Value = "'18"
$ws.Range("D7").Value = "'2.987"
$ws.Range("G7").Value = "'18"
$ws.Range("D8").Value = "'0.8102"
$ws.Range("G8").Value = "'18"
$ws.Range("D9").Value = "'0.8392"
$ws.Range("G9").Value = "'18"
$ws.Range("G10").Value = "'18"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03329"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").Value = "'18"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06942"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "'18"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02844"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "'18"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09410"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "'18"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001525"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("G15").Value = "'18"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0005978"
$ws.Range("E16").Value = "15OneONEWorstin24h"
$ws.Range("G16").Value = "'18"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006250"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("G17").Value = "'18"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.499"
$ws.Range("E18").Value = "17LEOLEO"
$ws.Range("G18").Value = "'18"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.091"
$ws.Range("E19").Value = "18BTSETokenBTSE"
$ws.Range("G19").Value = "'18"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3170"
$ws.Range("E20").Value = "19BitpandaEcosystemTokenBEST"
$ws.Range("G20").Value = "'18"
$ws.Range("D21").Value = "'0.1291"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'3.741"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'0.04693"
$ws.Range("G23").Value = "'18"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.001244"
$ws.Range("G25").Value = "'18"
$ws.Range("D26").Value = "'0.004524"
$ws.Range("G26").Value = "'18"
$ws.Range("G27").Value = "'18"
$ws.Range("D28").Value = "'0.0001939"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("G38").Value = "'18"
$ws.Range("G39").Value = "'18"
$ws.Range("G40").Value = "'18"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1052"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("G41").Value = "'18"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002722"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Value = "'18"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.006209"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("G43").Value = "'18"
$ws.Range("D44").Value = "'0.008383"
$ws.Range("G44").Value = "'18"
$ws.Range("D45").Value = "'0.00005258"
$ws.Range("G45").Value = "'18"
$ws.Range("G46").Value = "'18"
$ws.Range("D47").Value = "'0.2199"
$ws.Range("G47").Value = "'18"
$ws.Range("G48").Value = "'18"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("G49").Value = "'18"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("G50").Value = "'18"
$ws.Range("G51").Value = "'18"